$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Number1"
$ws.Range("B1").Value = "Number2"
$ws.Range("C1").Value = "Result"

$ws.Range("B2").Select()
